$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17 ("Rogue River National Forest") becomes "New Folder" (in-place rename,
# keeping the shared-string slot rather than appending a brand new one).
$ws.Cells.Item(17, 1).Value = "New Folder"

# Rows that are no longer needed are removed entirely (not just cleared),
# so later rows shift up to close the gaps. Deleting from the bottom up
# keeps the remaining row numbers stable while we work.
$rowsToDelete = @(25, 22, 21, 19, 16, 14, 12, 11, 10, 9, 7, 5, 2)

foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}
